$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that change per row: A, B, D, E, F, G, H, Q
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q")

# Capture current (pre-edit) values for rows 2, 3, 4 for the affected columns
$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("$col" + "2").Value()
    $row3[$col] = $ws.Range("$col" + "3").Value()
    $row4[$col] = $ws.Range("$col" + "4").Value()
}

# Cyclic shift: row3 -> row2, row4 -> row3, row2 -> row4
foreach ($col in $cols) {
    $ws.Range("$col" + "2").Value = $row3[$col]
    $ws.Range("$col" + "3").Value = $row4[$col]
    $ws.Range("$col" + "4").Value = $row2[$col]
}
